$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Row 58: replace the F58 "pend" marker with an E58 "hecho" marker.
$ws.Range("F58").ClearContents()
$ws.Range("E58").Value = "hecho"

# Row 74 previously held an empty (but styled) cell in B74; give it text and
# append 8 new backlog rows below it (75-82), matching the B-column task list.
$ws.Range("B74").Value = "investigar incrutacion de video"
$ws.Range("B75").Value = "acomodar apk"
$ws.Range("B76").Value = "cortar y listar videos"
$ws.Range("B77").Value = "agregar videos al sistema"
$ws.Range("B78").Value = "hacer nuevas capturas para la interfaz"
$ws.Range("B79").Value = "completar marco teorico con info de discapacitados"
$ws.Range("B80").Value = "terminar conclusion y  lineas futuras"
$ws.Range("B81").Value = "hacer manual de usuario"
$ws.Range("B82").Value = "revision total final de la documentacion"

# Carry the row-74 style (medium borders / teal fill / wrap-top) down across
# all the newly added rows so 75-82 look the same as the rest of the list.
$ws.Range("B75:B82").WrapText = $true
$ws.Range("B75:B82").VerticalAlignment = -4160
$ws.Range("B75:B82").Font.Bold = $true
$ws.Range("B75:B82").Font.Name = "Calibri"
$ws.Range("B75:B82").Font.Size = 11
$ws.Range("B75:B82").Borders.Item(8).LineStyle = 1
$ws.Range("B75:B82").Borders.Item(8).Weight = -4138
$ws.Range("B75:B82").Borders.Item(9).LineStyle = 1
$ws.Range("B75:B82").Borders.Item(9).Weight = -4138
$ws.Range("B75:B82").Interior.Pattern = -4124

# Reposition the view like the saved workbook (scrolled a bit further up,
# selection sitting on B73 instead of B74).
$ws.Application.ActiveWindow.ScrollRow = 63
$ws.Range("B73").Select()
